$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15999.25
$ws.Range("I18").Value = 15999.25
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 15999.25
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -15715.25
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 684.5714
$ws.Range("I33").Value = 407.63635
$ws.Range("K33").Value = 407.63635
$ws.Range("M33").Value = -178.63635
$ws.Range("H62").Value = 3042.5625
$ws.Range("I62").Value = 2797.2856
$ws.Range("K62").Value = 2797.2856
$ws.Range("M62").Value = -2173.2856
$ws.Range("H65").Value = 3042.5625
$ws.Range("I65").Value = 2797.2856
$ws.Range("K65").Value = 13986.428
$ws.Range("M65").Value = -10866.428
$ws.Range("H80").Value = 53633.344
$ws.Range("I80").Value = 84441.125
$ws.Range("J80").Value = 820
$ws.Range("K80").Value = 253323.375
$ws.Range("L80").Value = 2460
$ws.Range("M80").Value = -252325.375
$ws.Range("N80").Value = -4456
$ws.Range("H83").Value = 53633.344
$ws.Range("I83").Value = 84441.125
$ws.Range("J83").Value = 820
$ws.Range("K83").Value = 759970.125
$ws.Range("L83").Value = 7380
$ws.Range("M83").Value = -754978.125
$ws.Range("N83").Value = -17364
$ws.Range("H86").Value = 20609.854
$ws.Range("I86").Value = 3885.1875
$ws.Range("K86").Value = 3885.1875
$ws.Range("M86").Value = -2762.1875
$ws.Range("H89").Value = 20609.854
$ws.Range("I89").Value = 3885.1875
$ws.Range("K89").Value = 19425.9375
$ws.Range("M89").Value = -13809.9375
$ws.Range("H100").Value = 89663.47
$ws.Range("J100").Value = 51046
$ws.Range("L100").Value = 51046
$ws.Range("N100").Value = -52128
$ws.Range("H106").Value = 6179680
$ws.Range("I106").Value = 6865219.5
$ws.Range("J106").Value = 9825
$ws.Range("K106").Value = 6865219.5
$ws.Range("L106").Value = 9825
$ws.Range("M106").Value = -6864588.5
$ws.Range("N106").Value = -11087
$ws.Range("H132").Value = 1564934.6
$ws.Range("I132").Value = 2456.5344
$ws.Range("K132").Value = 7369.6032
$ws.Range("M132").Value = -4839.6032
$ws.Range("H137").Value = 7495.5
$ws.Range("I137").Value = 9486.484
$ws.Range("J137").Value = 2441.4614
$ws.Range("K137").Value = 28459.452
$ws.Range("L137").Value = 7324.3842
$ws.Range("M137").Value = -25909.452
$ws.Range("N137").Value = -12424.3842
$ws.Range("H138").Value = 296006.4
$ws.Range("J138").Value = 5005.115
$ws.Range("L138").Value = 15015.345
$ws.Range("N138").Value = -25295.345

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 188574.73
$ws.Range("I45").Value = 339551.34
$ws.Range("J45").Value = 7402.8
$ws.Range("K45").Value = 339551.34
$ws.Range("L45").Value = 7402.8
$ws.Range("M45").Value = -339174.34
$ws.Range("N45").Value = -8156.8
$ws.Range("H61").Value = 9773.852000000001
$ws.Range("I61").Value = 11993.5
$ws.Range("K61").Value = 11993.5
$ws.Range("M61").Value = -11781.5
$ws.Range("H74").Value = 6086.609
$ws.Range("I74").Value = 10740.1
$ws.Range("K74").Value = 10740.1
$ws.Range("M74").Value = -9866.1
$ws.Range("H77").Value = 6086.609
$ws.Range("I77").Value = 10740.1
$ws.Range("K77").Value = 53700.5
$ws.Range("M77").Value = -49332.5
$ws.Range("H122").Value = 942587.9
$ws.Range("I122").Value = 4652.8696
$ws.Range("K122").Value = 13958.6088
$ws.Range("M122").Value = -11508.6088
$ws.Range("H132").Value = 2822.9565
$ws.Range("I132").Value = 1233.125
$ws.Range("K132").Value = 3699.375
$ws.Range("M132").Value = -1169.375
$ws.Range("H133").Value = 70083
$ws.Range("J133").Value = 70083
$ws.Range("L133").Value = 70083
$ws.Range("N133").Value = -75143
$ws.Range("H136").Value = 9773.852000000001
$ws.Range("I136").Value = 11993.5
$ws.Range("K136").Value = 35980.5
$ws.Range("M136").Value = -33430.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5772.727
$ws.Range("I86").Value = 6058.3335
$ws.Range("J86").Value = 4487.5
$ws.Range("K86").Value = 6058.3335
$ws.Range("L86").Value = 4487.5
$ws.Range("M86").Value = -4935.3335
$ws.Range("N86").Value = -6733.5
$ws.Range("H89").Value = 5772.727
$ws.Range("I89").Value = 6058.3335
$ws.Range("J89").Value = 4487.5
$ws.Range("K89").Value = 30291.6675
$ws.Range("L89").Value = 22437.5
$ws.Range("M89").Value = -24675.6675
$ws.Range("N89").Value = -33669.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 601
$ws.Range("I22").Value = 601
$ws.Range("K22").Value = 601
$ws.Range("M22").Value = -251
$ws.Range("H31").Value = 9281.218999999999
$ws.Range("I31").Value = 13296.929
$ws.Range("K31").Value = 13296.929
$ws.Range("M31").Value = -13001.929
$ws.Range("H34").Value = 9281.218999999999
$ws.Range("I34").Value = 13296.929
$ws.Range("K34").Value = 13296.929
$ws.Range("M34").Value = -13094.929
$ws.Range("H107").Value = 40004964
$ws.Range("I107").Value = 47624890
$ws.Range("J107").Value = 363
$ws.Range("K107").Value = 47624890
$ws.Range("L107").Value = 363
$ws.Range("M107").Value = -47622970
$ws.Range("N107").Value = -4203
$ws.Range("H132").Value = 1165.7949
$ws.Range("I132").Value = 1095.9354
$ws.Range("K132").Value = 3287.8062
$ws.Range("M132").Value = -757.8062
$ws.Range("H141").Value = 181370.53
$ws.Range("J141").Value = 192891.62
$ws.Range("L141").Value = 192891.62
$ws.Range("N141").Value = -203251.62

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 30
$ws.Range("K3").Value = 90
$ws.Range("M3").Value = 22
$ws.Range("H7").Value = 89.75
$ws.Range("I7").Value = 83.35714
$ws.Range("K7").Value = 250.07142
$ws.Range("M7").Value = -138.07142
$ws.Range("H80").Value = 108005.234
$ws.Range("I80").Value = 192885.58
$ws.Range("J80").Value = 65565.07000000001
$ws.Range("K80").Value = 578656.74
$ws.Range("L80").Value = 196695.21
$ws.Range("M80").Value = -577720.74
$ws.Range("N80").Value = -198567.21
$ws.Range("H83").Value = 108005.234
$ws.Range("I83").Value = 192885.58
$ws.Range("J83").Value = 65565.07000000001
$ws.Range("K83").Value = 1735970.22
$ws.Range("L83").Value = 590085.6300000001
$ws.Range("M83").Value = -1731290.22
$ws.Range("N83").Value = -599445.6300000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 257207.2
$ws.Range("I20").Value = 2500000
$ws.Range("K20").Value = 2500000
$ws.Range("M20").Value = -2499755
$ws.Range("H24").Value = 1148577.1
$ws.Range("H45").Value = 30000
$ws.Range("J45").Value = 30000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -31118
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H132").Value = 3345.7632
$ws.Range("I132").Value = 3421.389
$ws.Range("K132").Value = 10264.167
$ws.Range("M132").Value = -7734.167000000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 499251.44
$ws.Range("I132").Value = 1242777
$ws.Range("J132").Value = 3567.7222
$ws.Range("K132").Value = 3728331
$ws.Range("L132").Value = 10703.1666
$ws.Range("M132").Value = -3725801
$ws.Range("N132").Value = -15763.1666

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 260540.4
$ws.Range("I62").Value = 636000.7
$ws.Range("K62").Value = 636000.7
$ws.Range("M62").Value = -635376.7
$ws.Range("H65").Value = 260540.4
$ws.Range("I65").Value = 636000.7
$ws.Range("K65").Value = 3180003.5
$ws.Range("M65").Value = -3176883.5
$ws.Range("H132").Value = 7816.277
$ws.Range("I132").Value = 9017.021000000001
$ws.Range("J132").Value = 4909.2104
$ws.Range("K132").Value = 27051.063
$ws.Range("L132").Value = 14727.6312
$ws.Range("M132").Value = -24521.063
$ws.Range("N132").Value = -19787.6312
